$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C ("Förändrad") date value from 45181 to 45182 for every
# existing data row (rows 2 through 524).
for ($i = 2; $i -le 524; $i++) {
    $ws.Range("C$i").Value = 45182
}

# Row 524 gains an explicit custom row height (ht="15" customHeight="1").
$ws.Rows.Item(524).RowHeight = 15

# Append new row 525: "A 42248-2023"
$ws.Range("A525").Value = "A 42248-2023"
$ws.Range("B525").Value = 45180
$ws.Range("B525").NumberFormat = "YYYY-MM-DD"
$ws.Range("C525").Value = 45182
$ws.Range("C525").NumberFormat = "YYYY-MM-DD"
$ws.Range("D525").Value = "SKÅNE LÄN"
$ws.Range("E525").Value = "KRISTIANSTAD"
$ws.Range("G525").Value = 6.1
$ws.Range("H525").Value = 0
$ws.Range("I525").Value = 0
$ws.Range("J525").Value = 0
$ws.Range("K525").Value = 0
$ws.Range("L525").Value = 0
$ws.Range("M525").Value = 0
$ws.Range("N525").Value = 0
$ws.Range("O525").Value = 0
$ws.Range("P525").Value = 0
$ws.Range("Q525").Value = 0
$ws.Range("R525").Value = $ws.Range("R524").Value
$ws.Range("R525").WrapText = $True
$ws.Rows.Item(525).RowHeight = 15

# Append new row 526: "A 42293-2023"
$ws.Range("A526").Value = "A 42293-2023"
$ws.Range("B526").Value = 45180
$ws.Range("B526").NumberFormat = "YYYY-MM-DD"
$ws.Range("C526").Value = 45182
$ws.Range("C526").NumberFormat = "YYYY-MM-DD"
$ws.Range("D526").Value = "SKÅNE LÄN"
$ws.Range("E526").Value = "KRISTIANSTAD"
$ws.Range("G526").Value = 0.6
$ws.Range("H526").Value = 0
$ws.Range("I526").Value = 0
$ws.Range("J526").Value = 0
$ws.Range("K526").Value = 0
$ws.Range("L526").Value = 0
$ws.Range("M526").Value = 0
$ws.Range("N526").Value = 0
$ws.Range("O526").Value = 0
$ws.Range("P526").Value = 0
$ws.Range("Q526").Value = 0
$ws.Range("R526").Value = $ws.Range("R524").Value
$ws.Range("R526").WrapText = $True
